{"js": "// Load all paragraphs in the document body so we can locate the two\n// target paragraphs by their current text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Change 1 --------------------------------------------------------\n// Paragraph \">>>  your stuff after this line >>>\" is currently split\n// across three runs (with proof-reading marks around the middle run).\n// Re-write it as a single run with the identical combined text, which\n// collapses the runs and drops the now-redundant proofErr markers.\nconst marker = \">  your\";\nconst quoteLinePara = paragraphs.items.find(\n  (p) => p.text.indexOf(marker) !== -1\n);\nif (quoteLinePara) {\n  quoteLinePara.insertText(\">>>  your stuff after this line >>>\", \"Replace\");\n}\n\n// --- Change 2 --------------------------------------------------------\n// Paragraph \"Baz changes\" (with a _GoBack bookmark sitting between the\n// two runs \"Baz chan\" and \"ges\") becomes \"David Fasullo changes for\n// assignment 1\". We replace each existing run's text in place (scoped\n// search within the paragraph) so the bookmark position is preserved.\nconst signaturePara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Baz chan\") !== -1\n);\nif (signaturePara) {\n  const firstRun = signaturePara.search(\"Baz chan\", { matchCase: true });\n  firstRun.load(\"items\");\n  await context.sync();\n  if (firstRun.items.length > 0) {\n    firstRun.items[0].insertText(\"David Fasullo\", \"Replace\");\n  }\n  await context.sync();\n\n  const secondRun = signaturePara.search(\"ges\", { matchCase: true });\n  secondRun.load(\"items\");\n  await context.sync();\n  if (secondRun.items.length > 0) {\n    secondRun.items[0].insertText(\" changes for assignment 1\", \"Replace\");\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two target paragraphs by content instead of a hard-coded\n# index, so the script is resilient to minor structural differences.\n$quoteLineIndex = $null\n$signatureLineIndex = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*>  your*\") {\n        $quoteLineIndex = $i\n    }\n    if ($t -like \"*Baz chan*\") {\n        $signatureLineIndex = $i\n    }\n}\n\n# --- Change 1 ----------------------------------------------------------\n# \">>>  your stuff after this line >>>\" is split across three runs\n# (with proofErr gramStart/gramEnd markers around the middle run).\n# Find/Replace the whole phrase with itself so Word collapses it into a\n# single run with identical text and drops the now-unneeded proofErr\n# markers.\nif ($quoteLineIndex) {\n    $quoteRange = $d.Paragraphs($quoteLineIndex).Range\n    $quoteFind = $quoteRange.Find\n    $quoteFind.ClearFormatting()\n    $quoteFind.Replacement.ClearFormatting()\n    $quoteText = \">>>  your stuff after this line >>>\"\n    $quoteFind.Execute($quoteText, $false, $false, $false, $false, $false, $true, 1, $false, $quoteText, 2) | Out-Null\n}\n\n# --- Change 2 ----------------------------------------------------------\n# \"Baz changes\" (with a _GoBack bookmark between the runs \"Baz chan\" and\n# \"ges\") becomes \"David Fasullo changes for assignment 1\". Replace each\n# run's text separately, scoped (via Paragraphs(...).Range) to this\n# paragraph only, so the bookmark's position between the two runs is\n# preserved and no other occurrence of \"ges\" in the document is touched.\nif ($signatureLineIndex) {\n    $firstRange = $d.Paragraphs($signatureLineIndex).Range\n    $firstFind = $firstRange.Find\n    $firstFind.ClearFormatting()\n    $firstFind.Replacement.ClearFormatting()\n    $firstFind.Execute(\"Baz chan\", $false, $false, $false, $false, $false, $true, 1, $false, \"David Fasullo\", 2) | Out-Null\n\n    $secondRange = $d.Paragraphs($signatureLineIndex).Range\n    $secondFind = $secondRange.Find\n    $secondFind.ClearFormatting()\n    $secondFind.Replacement.ClearFormatting()\n    $secondFind.Execute(\"ges\", $false, $false, $false, $false, $false, $true, 1, $false, \" changes for assignment 1\", 2) | Out-Null\n}\n"}
